$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titles = @(
  'Enrique Iglesias (feat. Nicole Scherzinger) - Heartbeat',
  'Sebastián Yatra - Adiós (Letra)',
  'Alan Walker - Lily (Lyrics) ft. K391, Emelie Hollow',
  'Taylor Swift - Out Of The Woods (Lyrics)',
  'One Direction - Story Of My Life',
  'Taylor Swift - Cruel Summer (Lyrics)',
  'Lewis Capaldi - Before You Go (Lyrics)',
  'Harry Styles - Sign of the Times',
  'Sunset Sons - VROL',
  'Passenger - Patient Love (Lyrics)',
  'Billie Eilish - lovely (Lyrics) ft. Khalid',
  'Arash feat.Helena - Angels Lullaby(Lyrics)',
  'Billie Eilish - Birds Of A Feather (Lyrics)',
  'Zivert – Goodbye | 2025',
  'Alphaville - Forever Young (Lyrics)',
  'Zivert - DEL MAR |  2021',
  'E.Satie - Gnossienne N.1 (Piano)',
  'Taylor Swift - Begin Again',
  'Lana Del Rey - Young and Beautiful',
  'Elyanna - Enta Eih (Hijazi Remix)',
  'La Tormenta De Arena - Dorian (letra)',
  'Lil Wayne feat. Bruno Mars - Mirror (Lyrics)',
  'Arash feat. Helena - Pure Love (Official Video)',
  '2CELLOS - Love Story',
  'Wiz Khalifa - See You Again ft. Charlie Puth',
  'twenty one pilots - Heathens',
  'the luka state - bring this all together',
  'New Moon - Edward leaves - Alexandre Desplat',
  'Каспийский Груз - Греет feat. Loc-Dog',
  'Justin Bieber - Baby',
  'OST Autumn in My Heart - Romance',
  'Sunset Sons - Remember',
  'Irakliy - Ya s toboy(cover)',
  'Morandi - Angels [Official Video]',
  'Passenger - Hell Or High Water',
  'Tones and I: Dance Monkey (US TV Debut)',
  'Skylar Grey - Love The Way You Lie',
  'Loving Calibri - Fed Up With Us',
  'Alexandre Desplat - New Moon (The Meadow)',
  'Каспийский Груз - С ней живой',
  'Passenger - All the little lights',
  'Justin Bieber & benny blanco - Lonely (Official Acoustic Video)',
  'Юлия Савичева – Москва-Владивосток',
  'For You - Liam Payne & Rita Ora',
  'Maksim - Vetrom stat (cover)',
  'Burito - Пока город спит',
  'Skylar Grey - I Know You (Lyrics)',
  'Elley Duhe - Middle Of The Night',
  'Rafet El Roman ft. Sinem - Seni Seviyorum',
  'Justin Bieber - Lonely (acoustic)',
  'HammAli & Navai - Птичка (cover)',
  'Broken Angel (Albert Vishi ft. Taulant Sllamniku Cover)',
  '서태지와 아이들   이 밤이 깊어가지만 (가사 첨부)',
  'Albert Vishi - My Time (Lyrics)',
  'ZAYN - Dusk Till Dawn (Official Video) ft. Sia',
  'Enrique Iglesias - Be With You (Lyrics)',
  'Tired Of Being Sorry (Laisse Le Destin L''Emporter)',
  'Sayyod feat. Rahnama - Tres metros sobre el cielo',
  'Passenger | Survivors (Official Video)',
  'Cody Francis - Rose In The Garden',
  'Coldplay - Hymn For The Weekend',
  'Enrique Iglesias - Ring My Bells (lyrics)',
  'Natalie Imbruglia - Torn',
  'Danna Paola, Sebastián Yatra - No Bailes Sola (Letra)',
  'Florida Georgia Line - Simple (Lyrics)',
  '2CELLOS - Fragile [LIVE at Arena Pula]',
  'Escape (Rosaline OST)',
  'Drake - God''s plan',
  'I Got Summer On My Mind (Still Dre Remix)',
  'Tim Odell - Another Love',
  'Passenger - To Be Free',
  'Alan Walker, Sabrina Carpenter &amp; Farruko  - On My Way',
  'Adele - Skyfall',
  'Lana Del Rey - Summertime Sadness (Official Music Video)',
  'Enrique Iglesias - Finally found you (lyrics) ft.Sammy Adams',
  'Yalın - Yeniden',
  'Linkin Park - Numb (cover)',
  'Passenger - Young As The Morning, Old As The Sea',
  'МакSим - Ветром стать',
  'Craig David - Insomnia (Lyrics)',
  'somewhere only we know (Gustixa &amp; Rhianne)',
  'Sunset Sons - Somewhere Maybe (Official Audio)',
  'Loving Caliber - We Were Dancing In The Dark',
  'Let me down slowly',
  'Akcent - I''m Sorry (lyrics)',
  '24kGoldn - Mood Remix (Lyrics) ft. Justin Bieber, J Balvin, Iann Dior',
  'Best classical music: Beethoven, Mozart, Schubert,Bach...🎶',
  'Somewhere Only We Know - Keane | Shania Yan Cover',
  'Alvaro Soler - Magia [Letra]',
  'Taylor Swift - Red (Taylor''s Version)',
  'Havana feat. Yaar &amp; Kaiia - Last Night (Lyrics)',
  'Dua Lipa - Levitating',
  'Jay Sean - Maybe | Lyrics',
  'Ramz - Barking',
  'Mariage D''Amour',
  'Edward Cullen - Bella''s Lullaby',
  'Cry (Acoustic) - Jamestown Story',
  'Akcent - That''s My Name (lyrics)',
  'Edward Maya, Vika Jigulina - Stereo love (Radio Edit) (Lyrics)',
  'Elyanna - Enta Eh(slowed)',
  'K-maro - Femme Like U (Paroles)',
  'Bruno Mars - Grenade',
  'Don`t Say Goodbye - Jamestown Story',
  'Баста & Zivert - неболей',
  'Rafet El Roman - Senden Sonra',
  'Enrique Iglesias - EL BAÑO (Letra) ft. Bad Bunny',
  'Zara Larsson – Dont Worry Bout Me',
  'Alvaro Soler - El Mismo Sol (letra)',
  'Elyanna - Enta Eih (lyrics)',
  'Ava Max - Who''s Laughing Now',
  'Charlie Puth - We Dont Talk Anymore',
  'Havana feat. Yaar & Kaiia - Big Love (Official Video)',
  'Баста - Выпускной',
  'Charlie Puth - Attention',
  'Maher Zain - Thank you Allah',
  'Passenger - Hard To Say Im Sorry',
  'Craig David - Walking Away [Lyrics] 🎵',
  'Let Her Go (ft Ed Sheeran)',
  'Eminem - Lose Yourself ⭐',
  'Sunset Sons - I Can`t Wait (Official Audio)',
  'Sunset Sons - The River',
  'Burito - По волнам',
  'post malone - rockstar (feat. 21 savage)',
  'Akcent - Stay With Me (lyrics)',
  'El mismo sol - Álvaro Soler & Jennifer López',
  'Sebastián Yatra, Álvaro Díaz - A Dónde Van (Letra)',
  'Ava Max - Alone',
  'Duncan Lawrence - Arcade',
  'Lewis Capaldi - Someone You Loved',
  'Rafet El Roman & Derya - Unuturum Elbet',
  'Halsey - Sorry',
  'The Weeknd - Blinding Lights',
  'Nelly ft Kelly Rowland - Dilemma',
  'Taylor Swift - Back To December',
  'Lana Del Rey - hope is a dangerous thing',
  'Passenger - Why Can''t I Change (Lyrics)',
  '2CELLOS - Shape Of My Heart [Live at Arena di Verona]',
  'Skylar Grey - Everything I Need',
  'Richard Clayderman - Fur Elise',
  'Chansons d’Amour en Français  💕',
  'Ava Max - Sweet But Psycho',
  'Alvaro Soler - Solo Para Ti',
  '2CELLOS - Now We Are Free - Gladiator [OFFICIAL VIDEO]',
  'Duncan Laurence feat. FLETCHER – Arcade',
  'Coldplay - Hunger Games | Atlas',
  'The Chainsmokers - Something Just Like This',
  'Juice Wrld - Lucid Dreams ',
  'Enrique Iglesias - Tired Of Being Sorry (Lyrics)',
  'Sture Zetterberg - body to body(가사/번역)',
  'Selena Gomez - Buscando Amor',
  'Elyanna - Youm Wara Youm',
  'Arash feat. Helena - Broken Angel (Official Video)',
  'Alec Benjamin ~ Outrunning Karma',
  'Marka - Ayol',
  'Anna Hamilton - Bad Liar (cover)',
  'Kaleida - Think (Lyrics) John Wick soundtrack',
  'Darren Hayes - Insatiable (Official Music Video)',
  'Taylor Swift - Champagne Problems',
  'Camila Cabello - Havana',
  'Stephen Sanchez - Until I Found You',
  'Jamestown Story -Ashamed',
  'Sebastián Yatra - Devuélveme el Corazón (Letra)',
  'Linkin Park - Castle of Glass',
  'HAUSER - Song from a Secret Garden',
  'Rafet El Roman & Derya - Özledim (Düet)',
  'Nightcore - Diamond Heart',
  'Enya - Only Time (Official 4K Music Video)',
  'Passenger | The Wrong Direction (Official Video)',
  'Justin Bieber - Let Me Love You',
  'Taylor Swift - State of Grace',
  'Sunset Song - On The Road (Lyrics)',
  'I`m Sorry - Jamestown Story',
  'Selena Gomez - Adiós',
  'Passenger - Survivors',
  'Akcent feat. Sandra N - Amor Gitana',
  'Sunset Sons - Loa (Official Audio)',
  'The Chainsmokers - Don''t Let Me Down',
  'Reamonn - Tonight',
  'Passenger - Simple Song (Lyrics)',
  'Gym Class Heroes: Stereo Hearts ft. Adam Levine',
  'Alec Benjamin - I Built a Friend',
  'Selena Gomez - Bad Liar',
  'Beth Thornton - Something You Don''t Know',
  'ZAYN, Zhavia Ward - A Whole New World',
  'Jamestown Story - Broken Summer',
  'Shakira - La La La World Cup 2014',
  'Alvaro Soler - La Cintura (Letra)',
  'Mina/Dracula - Please, don''t make me love you',
  'Passenger | Young As The Morning, Old As The Sea',
  'Supergirl - Anna Naklab ft. Allen Farben',
  'Reflex - В первый раз (cover)',
  'JEON SOMI DUMB DUMB Lyrics (전소미 DUMB DUMB 가사)',
  'Skylar Grey - Moving Mountains',
  'Loving Caliber - You Set My World On Fire',
  'Passenger - Heart is on fire',
  'Yanni - Can''t Wait (Sensuous Chill)',
  'Linkin Park - Numb (lyrics|rock)',
  'Harry Styles - As It Was',
  'Linkin Park - In the end (rmx)',
  'Reamonn - Supergirl',
  'Miley Cyrus - Flowers',
  'Maggie Lindemann - Pretty Girl',
  'Linkin Park - In The End (Mellen Gi &amp; Tommee Profitt Remix)',
  'Alan Walker - Play (Lyrics) ft. K-391, Tungevaag, Mangoo',
  'Lana Del Rey - Henry, come on (Lyrics)',
  'Richard Clayderman - Love Story',
  '2CELLOS - Fields Of Gold [Live at Arena di Verona]',
  'Today is a Good Day',
  'Lana Del Rey - Video Games',
  'Maroon 5 - Girls Like You ft. Cardi B (Official Music Video)',
  'could I have this kiss forever - Iglesias & Whitney',
  'Selena Gomez - Selfish love',
  'Yalın - Zalim (Official Video)',
  'Maroon 5, Wiz Khalifa – Payphone',
  'Dua Lipa - Break My Heart (Lyrics)',
  'Mariah Carey - My All [Lyrics]',
  'Taylor Swift - Style',
  'Maher Zain - Insha Allah',
  'ARASH feat Helena - ONE DAY (Official Video)',
  'Eminem - Mockingbird (Lyrics)',
  'Clandestina (Cover)',
  'HammAli & Navai - Птичка',
  'Maxim Fadeev - Googoosha',
  'Shakira - Hips Don''t Lie',
  'Loving Caliber - I Wish You Were Mine',
  'Rauf Faik - детство (Official audio)',
  'Janji - Heroes Tonight',
  'Shakira - Can`t Remember to Forget You (Lyrics) ft. Rihanna',
  'Merk & Kremont - Sad Story (Out Of Luck)',
  'Naomi Scott - Speechless (Lyrics)',
  'Maroon 5 - Memories',
  'Shakira - Chantaje (letra)',
  'Skylar Grey - Invisible',
  'Laura Pausini - It’s Not Goodbye',
  'Mika - Relax, Take it Easy',
  'Craig David - Rise & Fall ft. Sting (Official Video)',
  'Shawn Mendes - In My Blood',
  'Rafet El Roman - Kalbine Sürgün Feat. Ezo',
  'Zivert - Life (English Version)',
  'Alec Benjamin - Let Me Down Slowly (feat. Alessia Cara)',
  'Maher Zain - For The Rest Of My Life',
  'Christina Perri - A Thousand Years [Official Music Video]',
  'Passenger - Holes',
  'ARASH feat. Helena - DOOSET DARAM (Official Video)',
  'Hands Like Houses - Torn',
  'We''re All Runners - Craig Reever (Lyrics)',
  'Alan Walker Style , Adele - Set Fire To The Rain (Albert Vishi Remix)',
  'Akcent - Chimie Intre Noi',
  '올인(All In) - MV_처음 그날처럼 (2003)',
  'Alvaro Soler - Sofia',
  'Zara Larsson - Lush Life',
  'Sting - Shape of My Heart',
  'Zara Larsson - This Ones For You',
  'INNA - Oare (Rock The Roof @ Bucharest)',
  'New Moon OST - Dreamcatcher - Alexandre Desplat',
  'Lady Gaga, Bruno Mars - Die With A Smile',
  'Sunset Sons - Know My Name (Official Audio)',
  'Cry - Jamestown Story',
  'Shawn Mendes - It''ll Be Okay',
  'Passenger - Let Her Go',
  'falling in love with someone you can''t have (a playlist)',
  'Sasha Sloan - Lie'
)

$links = @(
  'https://www.youtube.com/watch?v=2hs0rtK9xwk',
  'https://www.youtube.com/watch?v=DULoaFTPB60',
  'https://www.youtube.com/watch?v=sZ5OUc7Ccwo',
  'https://www.youtube.com/watch?v=FEzj8K1cL6w',
  'https://youtu.be/W-TE_Ys4iwM?si=RViOxRuaXxdz3pmm',
  'https://www.youtube.com/watch?v=P8T1rUpVdXE',
  'https://www.youtube.com/watch?v=c-KFmp9MMmQ',
  'https://www.youtube.com/watch?v=qN4ooNx77u0',
  'https://www.youtube.com/watch?v=vMIgQ36zhAg',
  'https://www.youtube.com/watch?v=6LcKdxaSZVU',
  'https://www.youtube.com/watch?v=8VLXHyHRXjc',
  'https://www.youtube.com/watch?v=MeHCr0e-8vk',
  'https://www.youtube.com/watch?v=d5gf9dXbPi0',
  'https://www.youtube.com/watch?v=7vowLoOBk4Q',
  'https://www.youtube.com/watch?v=-ym0I1D8PfA',
  'https://www.youtube.com/watch?v=mH8l-7tRVnc',
  'https://www.youtube.com/watch?v=X3JLOenXGUc',
  'https://youtu.be/cMPEd8m79Hw?si=9zE5-51p0xGyEgSO',
  'https://www.youtube.com/watch?v=o_1aF54DO60',
  'https://www.youtube.com/watch?v=zFtsS4A8fc0',
  'https://www.youtube.com/watch?v=28W-KrHjmK8',
  'https://www.youtube.com/watch?v=97xukmZfiGU',
  'https://www.youtube.com/watch?v=KWKL8hvHR6E',
  'https://www.youtube.com/watch?v=UdHopftQD3A',
  'https://www.youtube.com/watch?v=RgKAFK5djSk',
  'https://youtu.be/UprcpdwuwCg?si=O6_fwxx8TOkfjIXi',
  'https://www.youtube.com/watch?v=OcJ5EgxsWBg',
  'https://www.youtube.com/watch?v=MS4Tf9mr44M',
  'https://www.youtube.com/watch?v=kSy7h4h9iC4',
  'https://www.youtube.com/watch?v=khOFw2f4bQY',
  'https://www.youtube.com/watch?v=y7roujOll3M',
  'https://www.youtube.com/watch?v=PH_P12XqY9Y',
  'https://youtu.be/3WmdZOF5bKk?si=LcXY8Gohxxx4cZSA',
  'https://www.youtube.com/watch?v=G0w8CEzMBrY',
  'https://youtu.be/zgDbp5C74sU?si=R8Q5HZq2vzhGL57g',
  'https://www.youtube.com/watch?v=4iQxG8ZjYO8',
  'https://youtu.be/h_-JFUci0BM?si=SHiuHs1NdIjpN0WP',
  'https://youtu.be/n1NTv6Y4pxs?si=76WA3JI0TGILBHm7',
  'https://www.youtube.com/watch?v=7kY5bQEU5gQ',
  'https://www.youtube.com/watch?v=a-Vf8T55gd8',
  'https://youtu.be/OkxVxox--Io?si=AE4wj_c_uqTWGrbB',
  'https://www.youtube.com/watch?v=Cu5hhxP_prE',
  'https://www.youtube.com/watch?v=s58SJ9pTXHE',
  'https://www.youtube.com/watch?v=ovZl1SjI2GY',
  'https://youtu.be/kkzEs0gdvZI?si=Z456wgKuJd0aE_PA',
  'https://www.youtube.com/watch?v=CPjVFf4Y_TQ',
  'https://www.youtube.com/watch?v=EBriiJpRGc8',
  'https://www.youtube.com/watch?v=KLTMCPzRO64',
  'https://www.youtube.com/watch?v=B3OcAOzFOCc',
  'https://youtu.be/Cu5hhxP_prE?si=VRZVlVcLWqk8Dasg',
  'https://www.youtube.com/watch?v=5hpSD-54Dtg',
  'https://www.youtube.com/watch?v=5miHGQVFJm0',
  'https://www.youtube.com/watch?v=__SXVP2GmvM',
  'https://www.youtube.com/watch?v=EOf5TP0kYHA',
  'https://www.youtube.com/watch?v=tt2k8PGm-TI',
  'https://www.youtube.com/watch?v=bZXnan-4GHo',
  'https://www.youtube.com/watch?v=ytnWeRME0aM',
  'https://www.youtube.com/watch?v=ZkV90Ea47wk',
  'https://www.youtube.com/watch?v=vN0gaXS8dQE',
  'https://www.youtube.com/watch?v=JO4-j1LfoQQ',
  'https://www.youtube.com/watch?v=YykjpeuMNEk',
  'https://www.youtube.com/watch?v=vhI_gDs_ZMg',
  'https://www.youtube.com/watch?v=PyKpNBAv1Bw',
  'https://www.youtube.com/watch?v=_ULT1lySBHk',
  'https://www.youtube.com/watch?v=TuTDc9d_9yI',
  'https://www.youtube.com/watch?v=q_ymIjOyzRQ',
  'https://www.youtube.com/watch?v=M9b_z-LKE14',
  'https://www.youtube.com/watch?v=ScfgOVJiu_I',
  'https://www.youtube.com/watch?v=89LOsf8pDhY',
  'https://youtu.be/Jkj36B1YuDU?si=Yku5tRPe7avRNr2R',
  'https://youtu.be/hNd5pILkpSw?si=qiwZxiuS0yeiuOPs',
  'https://www.youtube.com/watch?v=dhYOPzcsbGM',
  'https://youtu.be/DeumyOzKqgI?si=Cok0dR7byK6pN682',
  'https://www.youtube.com/watch?v=TdrL3QxjyVw',
  'https://www.youtube.com/watch?v=0d3eJ6OZoI4',
  'https://www.youtube.com/watch?v=iGut_MVMcUY',
  'https://youtu.be/gHp-OjLOG5A?si=0abUDswbKz6rhQeX',
  'https://www.youtube.com/watch?v=_3L0K5jXJv4',
  'https://www.youtube.com/watch?v=q8Il2rhe3MI',
  'https://www.youtube.com/watch?v=1en-4eT_iDY',
  'https://www.youtube.com/watch?v=92izkAK5OA0',
  'https://www.youtube.com/watch?v=SHapfmLyBp0',
  'https://www.youtube.com/watch?v=P-QYRUPDAQ8',
  'https://youtu.be/50VNCymT-Cs?si=sEwBTlJCeuqL9LTD',
  'https://www.youtube.com/watch?v=MfTmraRvihQ',
  'https://www.youtube.com/watch?v=f1J4dRTMy9A',
  'https://www.youtube.com/watch?v=DxnDcH2NS5c',
  'https://www.youtube.com/watch?v=kLKqeyx_HjY',
  'https://www.youtube.com/watch?v=sE1NoQRvFls',
  'https://www.youtube.com/watch?v=R_rUYuFtNO4',
  'https://www.youtube.com/watch?v=i-Yuf5-zTec',
  'https://www.youtube.com/watch?v=j2c3tR_qfiQ',
  'https://www.youtube.com/watch?v=XAR1HvssHdE',
  'https://www.youtube.com/watch?v=Q0QKUU95bVc',
  'https://www.youtube.com/watch?v=1ej1SI4BRv8',
  'https://www.youtube.com/watch?v=zQME-ChSwNM',
  'https://www.youtube.com/watch?v=Lg3WGLYEelU',
  'https://www.youtube.com/watch?v=PAOAfYUVPl0',
  'https://www.youtube.com/watch?v=y9Kqb2z9Lzs',
  'https://www.youtube.com/watch?v=GU9z22e7QX4',
  'https://www.youtube.com/watch?v=NM7FdfqNhm4',
  'https://youtu.be/4YrzJ9RZ9qY',
  'https://www.youtube.com/watch?v=sFEpzJo6Iuc',
  'https://www.youtube.com/watch?v=u0PrLVWDU4M',
  'https://www.youtube.com/watch?v=Z2g8NAg1bbI',
  'https://www.youtube.com/watch?v=8BbtBnnnvCM',
  'https://youtu.be/u_tzZd9kIWg?si=y-s2yCVh4U2JLsJJ',
  'https://www.youtube.com/watch?v=qPLX-Cv0aIs',
  'https://www.youtube.com/watch?v=YQjmLhETTEA',
  'https://youtu.be/4JYSgIiSZSA?si=3v9kDuzvYJvWaOsO',
  'https://youtu.be/bpFVJJBgtXY?si=L2NuwOWGhmdKacwg',
  'https://www.youtube.com/watch?v=aVFNJBqj5vU',
  'https://youtu.be/t1-yL-xvklc?si=YZ1rS5hZtleOFOy1',
  'https://youtu.be/Oz5JDtkf1as',
  'https://youtu.be/RBrdl0v_anc?si=cu3qNsVyUIIzZGvv',
  'https://youtu.be/XCmOdVia9DE?si=60M6i15UUakuL7DH',
  'https://www.youtube.com/watch?v=8AwamgSDpdA',
  'https://youtu.be/HTcL9WkB_wg?si=ILXw9EaPM4GJyx29',
  'https://www.youtube.com/watch?v=PwclAKuCJUk',
  'https://www.youtube.com/watch?v=JuiegvRQ8dI',
  'https://www.youtube.com/watch?v=MCyEm1fViZQ',
  'https://youtu.be/jqyJ4xW2gb0?si=VgrA4JKMWkeWDIA5',
  'https://www.youtube.com/watch?v=9lQP9-F8kIQ',
  'https://www.youtube.com/watch?v=SmONLwqQiZE',
  'https://www.youtube.com/watch?v=l0qJUOgS4Qw',
  'https://www.youtube.com/watch?v=Mrj56WSTfxI',
  'https://youtu.be/omvW1cI-3xg?si=zHiFadZaUUpddcgu',
  'https://youtu.be/Qau6mObfSGM?si=RsrcZ0VUCOHaEwE4',
  'https://youtu.be/zABLecsR5UE?si=k3rryaA0P3O8JBhY',
  'https://www.youtube.com/watch?v=ScZFzmN-8XY',
  'https://youtu.be/CPAoMCo7tNw?si=2rEiXXCn6UcySUVZ',
  'https://www.youtube.com/watch?v=XwxLwG2_Sxk',
  'https://www.youtube.com/watch?v=O-B_bMh1hi4',
  'https://youtu.be/QUwxKWT6m7U?si=LNPBWKl0DqXIfOP2',
  'https://www.youtube.com/watch?v=qq9zYxW_uNo',
  'https://www.youtube.com/watch?v=cBCEaYVOgSA',
  'https://www.youtube.com/watch?v=jx1-NP9_YIA',
  'https://www.youtube.com/watch?v=9bCp7j3nC30',
  'https://www.youtube.com/watch?v=El8TEqRZ7ik',
  'https://www.youtube.com/watch?v=sHE6aKis69U',
  'https://youtu.be/2KBFD0aoZy8',
  'https://youtu.be/5D_A4IBWSv4?si=pgNinSqUyLBks6po',
  'https://www.youtube.com/watch?v=74CYIdYoQ5w',
  'https://www.youtube.com/watch?v=308v08mFWWc',
  'https://youtu.be/Lh3TokLzzmw?si=I5CcdBNIEuwDZvVT',
  'https://youtu.be/FM7MFYoylVs?si=TrbAGj-JAUeEJ4bd',
  'https://www.youtube.com/watch?v=_fh64GbFSw4',
  'https://www.youtube.com/watch?v=TkzWwNiBtqE',
  'https://www.youtube.com/watch?v=bIywxOrMFvY',
  'https://youtu.be/2P6EExu3H5s?si=f2hv9y52VqxnVOmL',
  'https://www.youtube.com/watch?v=ceS-wnW9zSw',
  'https://www.youtube.com/watch?v=hbqoaJ8sKsQ',
  'https://www.youtube.com/watch?v=SqGrXzIPoQo',
  'https://www.youtube.com/watch?v=U8Y_jRXW2vI',
  'https://youtu.be/5jfz3q9Z0RY?si=OHvyb7AMtM_wtAXc',
  'https://www.youtube.com/watch?v=FVtFcbBfNYw',
  'https://www.youtube.com/watch?v=9u7hGkL57N8',
  'https://www.youtube.com/watch?v=wMpqCRF7TKg',
  'https://youtu.be/HCjNJDNzw8Y?si=QjZAi7GPIc4ParOQ',
  'https://youtu.be/oIKuyj2GQtY',
  'https://www.youtube.com/watch?v=XEyOV7rACOQ',
  'https://www.youtube.com/watch?v=w9SBPJxP3t0',
  'https://www.youtube.com/watch?v=PPkJeWPP2nM',
  'https://www.youtube.com/watch?v=n9tw4tIUEoM',
  'https://www.youtube.com/watch?v=JJ1fR1X4NYk',
  'https://youtu.be/bcHoBDw4G10?si=auASu-G_c9NkS48Z',
  'https://www.youtube.com/watch?v=7wfYIMyS_dI',
  'https://www.youtube.com/watch?v=VvRVu78IHHo',
  'https://youtu.be/SMs0GnYze34?si=T-UORWGqJCoitcOM',
  'https://www.youtube.com/watch?v=gr4cqcqnAN0',
  'https://www.youtube.com/watch?v=NsKZ-5EDqPA',
  'https://www.youtube.com/watch?v=YnCMnsPH6d0',
  'https://youtu.be/9H_368c2Hzw?si=UOBGyTGbUe_fISFW',
  'https://www.youtube.com/watch?v=vN0gaXS8dQE',
  'https://www.youtube.com/watch?v=mXSryKIbE7g',
  'https://www.youtube.com/watch?v=9tXWQy7mMsM',
  'https://youtu.be/Io0fBr1XBUA?si=SUp9MdCXlOU_Vf5s',
  'https://youtu.be/jtoncUzV6nA?si=yULSO1-MxnAVV13i',
  'https://www.youtube.com/watch?v=4bBs9Q7XDPA',
  'https://www.youtube.com/watch?v=T3E9Wjbq44E',
  'https://www.youtube.com/watch?v=7KEpSl7DvBg',
  'https://www.youtube.com/watch?v=NZKXkD6EgBk',
  'https://www.youtube.com/watch?v=RftohIbwlqg',
  'https://www.youtube.com/watch?v=rg_zwK_sSEY',
  'https://www.youtube.com/watch?v=HCZecLh5o4Q',
  'https://youtu.be/2igups6VdcA?si=N5uu5genirJuWXWC',
  'https://www.youtube.com/watch?v=32lkXvYNwpI',
  'https://www.youtube.com/watch?v=i3BV6sQtyps',
  'https://www.youtube.com/watch?v=_3L0K5jXJv4',
  'https://www.youtube.com/watch?v=swBR4QnO3yE',
  'https://www.youtube.com/watch?v=LnwtmsjlRiQ',
  'https://www.youtube.com/watch?v=TfAzTYzBvTo',
  'https://youtu.be/S_0r3hYg78o?si=Be6GShy7mgRcl9Ha',
  'https://www.youtube.com/watch?v=nQ7SQVXkWr8',
  'https://youtu.be/kBqqlW6-99M?si=kXaaJTqhA4PaY6Gd',
  'https://www.youtube.com/watch?v=9kardLhsFrk',
  'https://youtu.be/8P0vKLHbtMg?si=HhXMHjE8vD2yeC_B',
  'https://youtu.be/Qfm6nfz1QNQ?si=3mMjYFpALij7GELl',
  'https://youtu.be/WNeLUngb-Xg?si=V95nGOt0sMvhQG7c',
  'https://youtu.be/ctmS5XX67Ek?si=NGZGPw0bcpfZciyi',
  'https://www.youtube.com/watch?v=iawgB2CDCrw',
  'https://www.youtube.com/watch?v=WEFJnYMz0Ec',
  'https://www.youtube.com/watch?v=WNeLUngb-Xg',
  'https://www.youtube.com/watch?v=trnx5XT0cZs',
  'https://www.youtube.com/watch?v=wasFuIuPh5k',
  'https://www.youtube.com/watch?v=ITswHbJPHhQ',
  'https://www.youtube.com/watch?v=STBa_TmxgK4',
  'https://youtu.be/9L4EjJqrz0c?si=x97RAvAA9IELRZPW',
  'https://www.youtube.com/watch?v=cE6wxDqdOV0',
  'https://www.youtube.com/watch?v=aJOTlE1K90k',
  'https://www.youtube.com/watch?v=NTyOwChDYV0',
  'https://youtu.be/9gqAq6kq5Ek?si=Gro32XWDuPLWzyIv',
  'https://www.youtube.com/watch?v=kPM5VXhpCfA',
  'https://www.youtube.com/watch?v=bbdsIR4UHDg',
  'https://www.youtube.com/watch?v=jgh8owCuX78',
  'https://www.youtube.com/watch?v=o4che1p-M4M',
  'https://www.youtube.com/watch?v=XAVLUYDtCCw',
  'https://youtu.be/8xXJyFNfiy8?si=XkqgGm4hEyZoqJe1',
  'https://www.youtube.com/watch?v=rjBsQ9SygnE',
  'https://www.youtube.com/watch?v=mkE8WSS51mA',
  'https://www.youtube.com/watch?v=UTgnOwpafpw',
  'https://www.youtube.com/watch?v=7T82Z7BwYfM',
  'https://www.youtube.com/watch?v=gqOoJXttEec',
  'https://youtu.be/p3pEe6aAJ4k?si=bzrAEs7c-zSwqBUo',
  'https://www.youtube.com/watch?v=5j9FKszXLag',
  'https://www.youtube.com/watch?v=WJF5Z1WRcqw',
  'https://www.youtube.com/watch?v=074rfF4RJZc',
  'https://www.youtube.com/watch?v=i_XM3u1_jZQ',
  'https://www.youtube.com/watch?v=8GotXeCwUnc',
  'https://www.youtube.com/watch?v=BaSf-ddZxB8',
  'https://www.youtube.com/watch?v=SlPhMPnQ58k&pp=ygUPbWFyb29uIG1lbW9yaWVz',
  'https://youtu.be/J76eQJP3UIQ?si=juYKqG_UCEta8y19',
  'https://www.youtube.com/watch?v=NVVrT_wNw_Y',
  'https://www.youtube.com/watch?v=onYQkI8S1UY',
  'https://youtu.be/EVDYmBrl02Q?si=ODB07HFZCtTtg4F4',
  'https://www.youtube.com/watch?v=pU2ukeS2JTE',
  'https://youtu.be/36tggrpRoTI?si=CiCfVdO8Oepjt4Rs',
  'https://www.youtube.com/watch?v=7I3h7czMJeI',
  'https://www.youtube.com/watch?v=mTecGII7cFA',
  'https://www.youtube.com/watch?v=gEo8IrFbecM',
  'https://youtu.be/PHbZ9SXHJwA?si=_7a2Gaka2oPEWrCQ',
  'https://www.youtube.com/watch?v=rtOvBOTyX00',
  'https://youtu.be/DeFWClW7skQ?si=hkIGl-CTTw-FbnLz',
  'https://www.youtube.com/watch?v=-yQ8kxikSJQ',
  'https://www.youtube.com/watch?v=M58IJO7N32s',
  'https://www.youtube.com/watch?v=a0u9wnYRp9I',
  'https://www.youtube.com/watch?v=cY1_o8yrILc',
  'https://www.youtube.com/watch?v=OK3HTPs1ccI',
  'https://www.youtube.com/watch?v=ddD9G7KQzx0',
  'https://youtu.be/ftI_Lp7LAuU?si=aOFT5Ral2-A_2PxG',
  'https://www.youtube.com/watch?v=tD4HCZe-tew',
  'https://youtu.be/pm3rDbXbZRI?si=7TxDuViBxhHGeZoU',
  'https://youtu.be/MoHnffhBwqs?si=_FGX4ucMtOTcD2to',
  'https://www.youtube.com/watch?v=3UOIIaYeDRI',
  'https://www.youtube.com/watch?v=39Kvcgug2J0',
  'https://www.youtube.com/watch?v=zgaCZOQCpp8',
  'https://www.youtube.com/watch?v=orMwK0veDVQ',
  'https://www.youtube.com/watch?v=S10klmMrCkc',
  'https://youtu.be/KrgJp7Z1Hv8?si=MOyY5rZzP-7kcfhM',
  'https://www.youtube.com/watch?v=_BCtgSCulIU',
  'https://www.youtube.com/watch?v=_K57AlI62V4',
  'https://youtu.be/AzjTJpzfB8U?si=PHYxAGETm1P1opd0'
)

for ($i = 0; $i -lt $titles.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 2).Value = $titles[$i]
  $ws.Cells.Item($row, 3).Value = $links[$i]
}

Write-Output "done"